# Apply the edits described by the commit diff:
#  - Rename "Isabel" -> "Snuisters" and "BESTELBON_ISABEL" -> "BESTELBON_SNUISTERS"
#    in the "Details" sheet.
#  - Make "Details" the active sheet (instead of "Order"), and set the
#    selected cell on "Details" to F3.

$wb = $excel.ActiveWorkbook

$wsDetails = $wb.Worksheets.Item("Details")

# Update the Details sheet values
$wsDetails.Range("A2").Value = "Snuisters"
$wsDetails.Range("C2").Value = "BESTELBON_SNUISTERS"

# Select the Details sheet, and a cell within it, then activate it so it
# becomes the workbook's active (displayed) sheet.
$wsDetails.Select()
$wsDetails.Range("F3").Select()
$wsDetails.Activate()
